$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("s6",  "s6_IMG_3174.jpeg",  "meltpatch", "289",  "486",  "104", "52", "6",   "2"),
    @("s7",  "s7_IMG_3177.jpeg",  "meltpatch", "2513", "1794", "104", "52", "48",  "2"),
    @("s8",  "s8_IMG_3179.jpeg",  "meltpatch", "3210", "1456", "104", "52", "19",  "2"),
    @("s9",  "s9_IMG_3175.jpeg",  "meltpatch", "3174", "2241", "104", "52", "7",   "2"),
    @("s10", "s10_IMG_3178.jpeg", "meltpatch", "2826", "1369", "104", "52", "172", "2")
)

$startRow = 7
$endRow = $startRow + $data.Length - 1
# Columns D:I hold digit-only strings (e.g. "289"); force Text so Excel
# keeps them as strings instead of auto-converting to numbers (columns
# A:C are never ambiguous - "s6", "s6_IMG_3174.jpeg", "meltpatch" - and
# need no special formatting).
$ws.Range("D$($startRow):I$($endRow)").NumberFormat = "@"

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $values = $data[$i]
    for ($col = 1; $col -le $values.Length; $col++) {
        $ws.Cells.Item($row, $col).Value = $values[$col - 1]
    }
}
